$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Dictionary")

# Insert a new column before the "unit" column (H) to hold the new
# "data source" field, shifting unit / unit_uri / accepted_values right.
$ws.Range("H1").EntireColumn.Insert()

# Populate the new column's header and per-row values (written in this
# order so the shared-string table is built up the same way it was in
# the source edit).
$ws.Range("H1").Value = "data source"
$ws.Range("H3").Value = "Drone survey"
$ws.Range("H2").Value = "Calculated in this report"
$ws.Range("H4").Value = "WSP Database, NuSEDS"

# The column insert shifts the old "unit_uri" hyperlinks from I2/I3 to
# J2/J3, but leaves the <hyperlinks> anchors pointing at the old cells.
# Rebuild all hyperlinks on the sheet against their correct cells.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "http://example.org/spsr/stock_recruit_index")
$ws.Hyperlinks.Add($ws.Range("J2"), "http://vocab.nerc.ac.uk/collection/P06/current/UPMS/")
$ws.Hyperlinks.Add($ws.Range("E3"), "http://example.org/spsr/spawner_count")
$ws.Hyperlinks.Add($ws.Range("J3"), "http://vocab.nerc.ac.uk/collection/P06/current/NOPC/")
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("J2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("J3").Style = "Hyperlink"

# Widen the "method" column and size the new "data source" column.
$ws.Columns.Item(7).ColumnWidth = 61
$ws.Columns.Item(8).ColumnWidth = 13.1

# Leave the workbook focused on the Data Dictionary tab.
[void]$ws.Activate()
[void]$ws.Range("I23").Select()
